$wb = $excel.ActiveWorkbook

$wsTrain = $wb.Worksheets.Item("Train")
$wsVal = $wb.Worksheets.Item("Val")
$wsTest = $wb.Worksheets.Item("Test")

# Train sheet updates
$wsTrain.Range("B2").Value = 251.43
$wsTrain.Range("B3").Value = 114.1
$wsTrain.Range("B4").Value = 83.55
$wsTrain.Range("B5").Value = 68.02
$wsTrain.Range("B6").Value = 57.38
$wsTrain.Range("B7").Value = 49.04
$wsTrain.Range("B8").Value = 41.77
$wsTrain.Range("B9").Value = 35.51
$wsTrain.Range("B10").Value = 30.42
$wsTrain.Range("B11").Value = 26.46
$wsTrain.Range("B12").Value = 23.28
$wsTrain.Range("B13").Value = 20.75
$wsTrain.Range("B14").Value = 18.77
$wsTrain.Range("B15").Value = 17.23
$wsTrain.Range("B16").Value = 16.03
$wsTrain.Range("B17").Value = 14.95
$wsTrain.Range("B18").Value = 14.09
$wsTrain.Range("B19").Value = 13.38
$wsTrain.Range("B20").Value = 12.7
$wsTrain.Range("B21").Value = 12.07
$wsTrain.Range("B22").Value = 11.67
$wsTrain.Range("B23").Value = 11.23
$wsTrain.Range("B24").Value = 10.89
$wsTrain.Range("B25").Value = 10.51
$wsTrain.Range("B26").Value = 10.31
$wsTrain.Range("B27").Value = 9.949999999999999
$wsTrain.Range("B28").Value = 9.66
$wsTrain.Range("B29").Value = 9.48
$wsTrain.Range("B30").Value = 9.23
$wsTrain.Range("B31").Value = 9.029999999999999
$wsTrain.Range("B32").Value = 8.890000000000001
$wsTrain.Range("B34").Value = 8.56
$wsTrain.Range("B35").Value = 8.449999999999999
$wsTrain.Range("B36").Value = 8.25
$wsTrain.Range("B37").Value = 8.039999999999999
$wsTrain.Range("B38").Value = 8.01
$wsTrain.Range("B39").Value = 7.89
$wsTrain.Range("B40").Value = 7.75
$wsTrain.Range("B41").Value = 5.28

# Val sheet updates
$wsVal.Range("B2").Value = 186.16
$wsVal.Range("B3").Value = 96.17
$wsVal.Range("B4").Value = 78.86
$wsVal.Range("B5").Value = 65.72
$wsVal.Range("B6").Value = 53.86
$wsVal.Range("B7").Value = 46.37
$wsVal.Range("B8").Value = 40.28
$wsVal.Range("B9").Value = 34.59
$wsVal.Range("B10").Value = 30.43
$wsVal.Range("B11").Value = 27.74
$wsVal.Range("B12").Value = 24.92
$wsVal.Range("B13").Value = 22.58
$wsVal.Range("B14").Value = 21.06
$wsVal.Range("B15").Value = 19.66
$wsVal.Range("B16").Value = 18.11
$wsVal.Range("B17").Value = 16.77
$wsVal.Range("B18").Value = 16.07
$wsVal.Range("B19").Value = 15.24
$wsVal.Range("B20").Value = 14.38
$wsVal.Range("B21").Value = 13.62
$wsVal.Range("B22").Value = 12.85
$wsVal.Range("B23").Value = 12.14
$wsVal.Range("B24").Value = 11.6
$wsVal.Range("B25").Value = 11.17
$wsVal.Range("B26").Value = 10.87
$wsVal.Range("B27").Value = 10.68
$wsVal.Range("B28").Value = 10.25
$wsVal.Range("B29").Value = 9.93
$wsVal.Range("B30").Value = 9.58
$wsVal.Range("B31").Value = 9.31
$wsVal.Range("B32").Value = 9.130000000000001
$wsVal.Range("B33").Value = 8.779999999999999
$wsVal.Range("B34").Value = 8.640000000000001
$wsVal.Range("B35").Value = 8.41
$wsVal.Range("B36").Value = 8.24
$wsVal.Range("B37").Value = 8.09
$wsVal.Range("B38").Value = 7.81
$wsVal.Range("B39").Value = 7.69
$wsVal.Range("B40").Value = 7.71
$wsVal.Range("B41").Value = 5.03

# Test sheet updates
$wsTest.Range("B2").Value = 5.29
